# Auto-generated edit script: updates the crypto price/volume table on Sheet1
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text does not look like a plain number (prices with multiple
#     thousands-separator dots, percent strings, coin names, URLs, unicode subscript
#     digits, etc.) -- a plain .Value assignment keeps these as text cells, matching
#     the original inlineStr cell type, since Excel only auto-converts to Number when
#     the text actually parses as one. ---
$ws.Range('D2').Value = '51.573.80'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '2.976.21'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E5').Value = '  +2.67%  '
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = '3.450.46'
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '2.984.95'
$ws.Range('E16').Value = '  +2.40%  '
$ws.Range('E17').Value = '  +7.25%  '
$ws.Range('D18').Value = '51.555.60'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('E25').Value = '  +8.12%  '
$ws.Range('E26').Value = '  +17.22%  '
$ws.Range('E27').Value = '  +19.33%  '
$ws.Range('E28').Value = '  +13.28%  '
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  -0.61%  '
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('E36').Value = '  +6.34%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('E45').Value = '  +14.72%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E48').Value = '  +4.53%  '
$ws.Range('D49').Value = '2.031.92'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('E50').Value = '  +6.06%  '
$ws.Range('E51').Value = '  +1.79%  '

# --- Cells whose new text IS a plain number (e.g. "104.31"): a bare .Value assignment
#     would make Excel silently re-type the cell as a Number, which would not match the
#     original text cell. Instead, write a text-literal formula (="104.31"), then Copy +
#     PasteSpecial(xlPasteValues=-4163) it onto itself -- that flattens the formula back
#     to a plain value while preserving the pasted values Text type (no re-parse as
#     Number) and without touching the cells style. Done one cell at a time: a single
#     multi-area Range("A1,A2,...") Copy/PasteSpecial flattens the areas sequentially and
#     would shift values into the wrong cells across the gaps.
$ws.Range('D5').Formula = '="383.05"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="104.31"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('D9').Formula = '="0.594"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('D10').Formula = '="37.09"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('D12').Formula = '="0.0847"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="18.31"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('D15').Formula = '="7.59"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('D17').Formula = '="0.998"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="3.27"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('D20').Formula = '="7.42"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="12.86"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('D23').Formula = '="69.14"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="263.55"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="2.92"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('D26').Formula = '="8.41"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('D27').Formula = '="7.75"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="1.00"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="25.97"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="9.87"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="34.65"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('D34').Formula = '="50.96"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="0.0451"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('D38').Formula = '="3.02"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('D39').Formula = '="16.96"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="2.59"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D42').Formula = '="1.83"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="122.28"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="21.71"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="0.278"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="3.31"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="0.0332"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('D51').Formula = '="5.14"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)

$excel.CutCopyMode = 0

Write-Output "Updated $($ws.Name): 92 cells across rows 2-51 (prices/volumes refreshed, Dai/EthereumClassic rows swapped)."
